$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Globo"
$ws.Range("B3").Value = "Bom Dia Inter"
$ws.Range("C3").Value = "Transporte"
$ws.Range("D3").Value = "2025-04-01T08:00"
$ws.Range("E3").Value = "Negativo"
$ws.Range("F3").Value = "Sem vans em Campos. Permissionários de vans do setor C fazem paralisação agora de manhã. Repórter *ao vivo* do ponto das vans do setor C, no Centro. Nenhuma van no setor C no local. Eles pedem que os repasses feitos pela prefeitura sejam realizados. Setor C atende Santo Eduardo, Santa maria e Morro do Coco. Por enquanto, só está funcionando para Travessão. Entrevista com passageiros. Eles alegam que foram pegos de surpresa e que estão pagando R$ 20 para conseguir voltar para casa. Em depoimento, passageira diz que ganha R$ 3,50 do patrão e tem que pagar R$ 20. Essa é a segunda paralisação deste ano. "
